$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 215; this shifts the existing rows 215-225
# down to 217-227 (new dimension becomes A1:T227) and inherits the row
# formatting (e.g. the date-format style on column D) from the rows below.
$ws.Rows.Item(215).Resize(2).Insert()

# New row 215: Murcott / Especial
$ws.Cells.Item(215, 1).Value = 5
$ws.Cells.Item(215, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(215, 3).Value = "Maule"
$ws.Cells.Item(215, 4).Value = 44509
$ws.Cells.Item(215, 5).Value = 7
$ws.Cells.Item(215, 6).Value = "Fruta"
$ws.Cells.Item(215, 7).Value = 100102
$ws.Cells.Item(215, 8).Value = "Cítricos"
$ws.Cells.Item(215, 9).Value = 100102004
$ws.Cells.Item(215, 10).Value = "Mandarina"
$ws.Cells.Item(215, 11).Value = "Murcott"
$ws.Cells.Item(215, 12).Value = "Especial"
$ws.Cells.Item(215, 13).Value = 100
$ws.Cells.Item(215, 14).Value = 7000
$ws.Cells.Item(215, 15).Value = 7000
$ws.Cells.Item(215, 16).Value = 7000
$ws.Cells.Item(215, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(215, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(215, 19).Value = 389
$ws.Cells.Item(215, 20).Value = 18

# New row 216: Murcott / Tercera
$ws.Cells.Item(216, 1).Value = 5
$ws.Cells.Item(216, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(216, 3).Value = "Maule"
$ws.Cells.Item(216, 4).Value = 44509
$ws.Cells.Item(216, 5).Value = 7
$ws.Cells.Item(216, 6).Value = "Fruta"
$ws.Cells.Item(216, 7).Value = 100102
$ws.Cells.Item(216, 8).Value = "Cítricos"
$ws.Cells.Item(216, 9).Value = 100102004
$ws.Cells.Item(216, 10).Value = "Mandarina"
$ws.Cells.Item(216, 11).Value = "Murcott"
$ws.Cells.Item(216, 12).Value = "Tercera"
$ws.Cells.Item(216, 13).Value = 120
$ws.Cells.Item(216, 14).Value = 4000
$ws.Cells.Item(216, 15).Value = 4000
$ws.Cells.Item(216, 16).Value = 4000
$ws.Cells.Item(216, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(216, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(216, 19).Value = 222
$ws.Cells.Item(216, 20).Value = 18
